$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "184"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "426816.00"

$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "995"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3151764.33"

$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "410"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1661698.25"

$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "114"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "540628.09"

$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "29"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "186643.82"

$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "40"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "80000.00"

$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "101"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "281752.38"

$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "425"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1338179.15"

$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "152"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "652202.10"

$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "47"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "241045.00"

$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "18"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "38621.00"

$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = "558"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1804484.47"

$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = "36"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "94105.00"

$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "162"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "423905.00"

$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "82"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "284298.00"

$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "23"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "91595.14"

$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "12"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "55500.00"

$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "53"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "135183.00"

$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "100"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "281768.17"

$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "579"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1988705.52"

$ws.Range("C52").NumberFormat = "@"
$ws.Range("C52").Value = "261"
$ws.Range("D52").NumberFormat = "@"
$ws.Range("D52").Value = "1134878.76"

$ws.Range("C54").NumberFormat = "@"
$ws.Range("C54").Value = "25"
$ws.Range("D54").NumberFormat = "@"
$ws.Range("D54").Value = "138213.00"

$ws.Range("C55").NumberFormat = "@"
$ws.Range("C55").Value = "20"
$ws.Range("D55").NumberFormat = "@"
$ws.Range("D55").Value = "62220.65"

$ws.Range("C80").NumberFormat = "@"
$ws.Range("C80").Value = "881"
$ws.Range("D80").NumberFormat = "@"
$ws.Range("D80").Value = "2810944.67"

$ws.Range("C81").NumberFormat = "@"
$ws.Range("C81").Value = "335"
$ws.Range("D81").NumberFormat = "@"
$ws.Range("D81").Value = "1351403.79"

$ws.Range("C84").NumberFormat = "@"
$ws.Range("C84").Value = "31"
$ws.Range("D84").NumberFormat = "@"
$ws.Range("D84").Value = "63500.00"
